$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Notes" sheet: the Funded? note used to be the single string "Y/N" in B13.
# Split it into two columns: B13 = "Y", C13 = "N".
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("B13").Value = "Y"
$notes.Range("C13").Value = "N"
$notes.Range("B14").Select()

# ---------------------------------------------------------------------------
# "Data" sheet: add real date / number formatting + validation so the sheet
# can grow (nrows no longer fixed to the 6 pre-formatted blank rows) and the
# Funded? column offers a Y/N pick list sourced from the Notes sheet.
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# The old sheet pre-formatted rows 2-7 in columns H/I as blank placeholders;
# with validation now covering the whole column there's no need for them, so
# clear both their contents and formatting first (this also shrinks the used
# range back down to just the header row).
$data.Range("H2:I7").Clear()

# Apply the date number format to the whole Begin/End/Submit Date columns
# (header row included).
$dateFormat = "m/d/yy;@"
$data.Range("H1:H1048576").NumberFormat = $dateFormat
$data.Range("I1:I1048576").NumberFormat = $dateFormat
$data.Range("J1:J1048576").NumberFormat = $dateFormat

# Give Begin Date / End Date columns a bit more room to show the new format.
$data.Columns.Item(8).ColumnWidth = 16.832
$data.Columns.Item(9).ColumnWidth = 7.998697916666667

# Decimal validation on the money columns.
$dAmt = $data.Range("D1:D1048576")
$dAmt.Validation.Delete()
$dAmt.Validation.Add(2, 1, 1, "-1000000", "1000000000")
$dAmt.Validation.IgnoreBlank = $false

$eAmt = $data.Range("E1:E1048576")
$eAmt.Validation.Delete()
$eAmt.Validation.Add(2, 1, 1, "-1000000", "1000000000")
$eAmt.Validation.IgnoreBlank = $false

# Date-picker validation on Begin Date / End Date / Submit Date.
$hCol = $data.Range("H1:H1048576")
$hCol.Validation.Delete()
$hCol.Validation.Add(4, 1, 1, "1", "73051")
$hCol.Validation.IgnoreBlank = $false

$iCol = $data.Range("I1:I1048576")
$iCol.Validation.Delete()
$iCol.Validation.Add(4, 1, 1, "1", "73051")

$jCol = $data.Range("J1:J1048576")
$jCol.Validation.Delete()
$jCol.Validation.Add(4, 1, 1, "1", "73051")

# Funded? column gets a Y/N pick list sourced from the Notes sheet (where the
# single "Y/N" note was just split into separate Y / N cells).
$fCol = $data.Range("F1:F1048576")
$fCol.Validation.Delete()
$fCol.Validation.Add(3, 1, 1, "=Notes!`$B`$13:`$C`$13")

$data.Range("G8").Select()
